$d = $word.ActiveDocument

# 1) "${ Clinical Status|radio|Urgent}" -> "${Clinical Status|radio|Urgent}"
#    (drop the stray space that was typed right after the opening "${")
$r1 = $d.Content
$null = $r1.Find.Execute("`${ Clinical Status", $true, $false, $false, $false, $false, $true, 1, $false, "`${Clinical Status", 2)

# 2) "Relative: ${title|Contraindications – relative}" -> same text, runs
#    re-typed/cleaned up so the "$", "{" and "title|...}" pieces aren't
#    split across so many separate runs any more.
$dash = [char]0x2013
$r2 = $d.Content
$findText = " `$`{title|Contraindications $dash relative}"
$null = $r2.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $findText, 2)
